$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 600.63336
$ws.Range("I41").Value = 503
$ws.Range("J41").Value = 698.26666
$ws.Range("K41").Value = 503
$ws.Range("L41").Value = 698.26666
$ws.Range("M41").Value = -63
$ws.Range("N41").Value = -1578.26666
$ws.Range("H43").Value = 12883.333
$ws.Range("I43").Value = 13500
$ws.Range("K43").Value = 13500
$ws.Range("M43").Value = -13431
$ws.Range("H46").Value = 2001994.6
$ws.Range("J46").Value = 2502043.2
$ws.Range("L46").Value = 7506129.600000001
$ws.Range("N46").Value = -7506367.600000001
$ws.Range("H60").Value = 2001994.6
$ws.Range("J60").Value = 2502043.2
$ws.Range("L60").Value = 7506129.600000001
$ws.Range("N60").Value = -7507097.600000001
$ws.Range("H82").Value = 55
$ws.Range("I82").Value = 55
$ws.Range("K82").Value = 165
$ws.Range("M82").Value = 241
$ws.Range("H85").Value = 55
$ws.Range("I85").Value = 55
$ws.Range("K85").Value = 165
$ws.Range("M85").Value = 1239
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H116").Value = 4916.375
$ws.Range("I116").Value = 4738.2
$ws.Range("J116").Value = 4997.364
$ws.Range("K116").Value = 4738.2
$ws.Range("L116").Value = 4997.364
$ws.Range("M116").Value = -1296.2
$ws.Range("N116").Value = -11881.364
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35520.38
$ws.Range("I32").Value = 20401.385
$ws.Range("J32").Value = 122874.555
$ws.Range("K32").Value = 20401.385
$ws.Range("L32").Value = 122874.555
$ws.Range("M32").Value = -20114.385
$ws.Range("N32").Value = -123448.555
$ws.Range("H61").Value = 1976.125
$ws.Range("I61").Value = 1555.6666
$ws.Range("K61").Value = 1555.6666
$ws.Range("M61").Value = -1343.6666
$ws.Range("H110").Value = 2358.3845
$ws.Range("I110").Value = 2471.182
$ws.Range("K110").Value = 2471.182
$ws.Range("M110").Value = -426.1819999999998
$ws.Range("H132").Value = 11647.091
$ws.Range("J132").Value = 5072
$ws.Range("L132").Value = 15216
$ws.Range("N132").Value = -20276
$ws.Range("H136").Value = 1976.125
$ws.Range("I136").Value = 1555.6666
$ws.Range("K136").Value = 4666.9998
$ws.Range("M136").Value = -2116.9998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18834.889
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 18834.889
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1217.875
$ws.Range("I16").Value = 1425.75
$ws.Range("J16").Value = 1010
$ws.Range("K16").Value = 1425.75
$ws.Range("L16").Value = 1010
$ws.Range("M16").Value = -1138.75
$ws.Range("N16").Value = -1584
$ws.Range("H22").Value = 425
$ws.Range("I22").Value = 425
$ws.Range("K22").Value = 425
$ws.Range("M22").Value = -75
$ws.Range("H58").Value = 7910.875
$ws.Range("I58").Value = 15369.714
$ws.Range("K58").Value = 15369.714
$ws.Range("M58").Value = -15166.714
$ws.Range("H59").Value = 13181.546
$ws.Range("J59").Value = 13999.8
$ws.Range("L59").Value = 13999.8
$ws.Range("N59").Value = -16289.8
$ws.Range("H74").Value = 67000
$ws.Range("J74").Value = 67000
$ws.Range("L74").Value = 67000
$ws.Range("N74").Value = -68748
$ws.Range("H77").Value = 67000
$ws.Range("J77").Value = 67000
$ws.Range("L77").Value = 201000
$ws.Range("N77").Value = -209736
$ws.Range("H86").Value = 6599.1113
$ws.Range("I86").Value = 5402.1665
$ws.Range("K86").Value = 5402.1665
$ws.Range("M86").Value = -4279.1665
$ws.Range("H89").Value = 6599.1113
$ws.Range("I89").Value = 5402.1665
$ws.Range("K89").Value = 27010.8325
$ws.Range("M89").Value = -21394.8325
$ws.Range("H99").Value = 168384.5
$ws.Range("J99").Value = 2314
$ws.Range("L99").Value = 2314
$ws.Range("N99").Value = -5310
$ws.Range("H113").Value = 1217.875
$ws.Range("I113").Value = 1425.75
$ws.Range("J113").Value = 1010
$ws.Range("K113").Value = 1425.75
$ws.Range("L113").Value = 1010
$ws.Range("M113").Value = 744.25
$ws.Range("N113").Value = -5350
$ws.Range("H126").Value = 168384.5
$ws.Range("J126").Value = 2314
$ws.Range("L126").Value = 6942
$ws.Range("N126").Value = -11882
$ws.Range("H132").Value = 3503
$ws.Range("I132").Value = 3315.2
$ws.Range("K132").Value = 9945.599999999999
$ws.Range("M132").Value = -7415.599999999999
$ws.Range("H136").Value = 7910.875
$ws.Range("I136").Value = 15369.714
$ws.Range("K136").Value = 46109.142
$ws.Range("M136").Value = -43559.142

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 129.18182
$ws.Range("J2").Value = 129.18182
$ws.Range("L2").Value = 775.0909199999999
$ws.Range("N2").Value = -1001.09092
$ws.Range("H93").Value = 4999
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H107").Value = 2155.7144
$ws.Range("I107").Value = 711.2727
$ws.Range("J107").Value = 3744.6
$ws.Range("K107").Value = 711.2727
$ws.Range("L107").Value = 3744.6
$ws.Range("M107").Value = 1208.7273
$ws.Range("N107").Value = -7584.6
$ws.Range("H122").Value = 1166
$ws.Range("I122").Value = 1260.6
$ws.Range("J122").Value = 1008.3333
$ws.Range("K122").Value = 3781.8
$ws.Range("L122").Value = 3024.9999
$ws.Range("M122").Value = -1331.8
$ws.Range("N122").Value = -7924.9999
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 4001.077
$ws.Range("I132").Value = 4001.077
$ws.Range("K132").Value = 12003.231
$ws.Range("M132").Value = -9473.231

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 37886.875
$ws.Range("J20").Value = 42127.855
$ws.Range("L20").Value = 42127.855
$ws.Range("N20").Value = -42579.855
$ws.Range("H30").Value = 4000
$ws.Range("I30").Value = 3500
$ws.Range("K30").Value = 3500
$ws.Range("M30").Value = -3392
$ws.Range("H132").Value = 7139.16
$ws.Range("J132").Value = 3904.9333
$ws.Range("L132").Value = 11714.7999
$ws.Range("N132").Value = -16774.7999
$ws.Range("H136").Value = 3505.9285
$ws.Range("I136").Value = 3017.6
$ws.Range("K136").Value = 9052.799999999999
$ws.Range("M136").Value = -6502.799999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15956.2
$ws.Range("J41").Value = 15935.75
$ws.Range("L41").Value = 15935.75
$ws.Range("N41").Value = -16715.75
$ws.Range("H54").Value = 39999.332
$ws.Range("I54").Value = 20000
$ws.Range("J54").Value = 49999
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 49999
$ws.Range("M54").Value = -19480
$ws.Range("N54").Value = -51039
$ws.Range("H70").Value = 20105
$ws.Range("J70").Value = 20105
$ws.Range("L70").Value = 20105
$ws.Range("N70").Value = -20735
$ws.Range("H73").Value = 20105
$ws.Range("J73").Value = 20105
$ws.Range("L73").Value = 20105
$ws.Range("N73").Value = -22289
$ws.Range("H96").Value = 1800
$ws.Range("J96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("N96").Value = -3746
$ws.Range("H113").Value = 714.7368
$ws.Range("I113").Value = 430.8
$ws.Range("J113").Value = 816.1429000000001
$ws.Range("K113").Value = 1292.4
$ws.Range("L113").Value = 2448.4287
$ws.Range("M113").Value = 877.5999999999999
$ws.Range("N113").Value = -6788.4287
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
